$d = $word.ActiveDocument

$replacements = @(
    @("2025-08-30 Saturday", "2025-08-31 Sunday"),
    @("26×55=1430", "54×42=2268"),
    @("29×40=1160", "28×91=2548"),
    @("14×17=238", "74×52=3848"),
    @("20×95=1900", "52×36=1872"),
    @("94×56=5264", "97×79=7663"),
    @("70×57=3990", "83×50=4150"),
    @("53×82=4346", "42×46=1932"),
    @("76×50=3800", "98×80=7840"),
    @("39×68=2652", "85×86=7310"),
    @("32×37=1184", "63×59=3717"),
    @("95×78=7410", "82×22=1804"),
    @("30×97=2910", "76×90=6840"),
    @("14×22=308", "28×53=1484"),
    @("85×45=3825", "24×56=1344"),
    @("92×40=3680", "47×22=1034"),
    @("44×91=4004", "16×98=1568"),
    @("30×16=480", "51×32=1632"),
    @("32×92=2944", "61×84=5124"),
    @("99×60=5940", "45×60=2700"),
    @("34×23=782", "86×20=1720"),
    @("70×92=6440", "40×46=1840"),
    @("65×37=2405", "51×29=1479"),
    @("23×14=322", "42×92=3864"),
    @("43×67=2881", "91×57=5187"),
    @("41×69=2829", "47×68=3196")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
